$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- New hyperparameter description texts (sheet2 D7 / D8) ---
$svrText = @"
Cost (Regularization Parameter): it controls the trade-off between having a smooth decision boundary and classifying the training points correctly. A smaller C encourages a wider margin, potentially allowing more training points to be misclassified but promoting a simpler model. A larger C penalizes misclassifications more, resulting in a more complex model. 
Kernel function: The kernel function determines the type of funtion to be applied in the covariates. Common choices include linear, polynomial, radial basis function (RBF), and sigmoid kernels. The choice of kernel depends genuinely on the characteristics of the data.
Gamma (Kernel Coefficient for RBF): It defines how far the influence of a single training example reaches. A small gamma will create a more generalized decision boundary, meaning a far reach, while a large gamma will create a more intricate boundary that may fit the training data more closely.
Epsilon: this hyperparameter is used to control the width of the margin. It defines a margin of tolerance where no penalty is given to errors, allowing some flexibility in fitting the data.
"@

$rfText = @"
ntree: number of trees in the forest. Increasing the number of trees generally improves performance, but it comes at the cost of increased computational complexity. 
mtry: The number of features to consider when looking for the best split at each node. It can be an integer (representing the exact number of features) or a float (representing a percentage of features).
Smaller values can reduce overfitting, while larger values may capture more information from the data.
maxnodes: The minimum number of samples required to split an internal node. Increasing this value can lead to a more conservative model, preventing splits that only capture noise. 
"@

$ws2.Range("D7").Value = $svrText
$ws2.Range("D8").Value = $rfText

# --- Row heights (sheet2, rows 7 & 8 grow to fit the new long text) ---
$ws2.Rows.Item(7).RowHeight = 316.8
$ws2.Rows.Item(8).RowHeight = 201.6

# --- Reposition the picture on sheet2 (no visual change) so the engine
# re-derives its "to" cell anchor against the new row heights; must be
# done before the column-width change below so the "from" anchor (which
# is based on columns A:G, unaffected by column D) is left untouched. ---
$pic = $ws2.Shapes.Item(1)
$pic.Left = 597.7048828125
$pic.Top = 0
$pic.Width = 401.0909251968503
$pic.Height = 432.33771653543306

# --- Column width (sheet2, column D widened for long text) ---
$ws2.Columns.Item(4).ColumnWidth = 54.25

# --- Selections on each sheet ---
$ws1.Range("D10").Select()

$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("D9").Select()

# best-effort: scroll sheet1's view too (topLeftCell A4)
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

# make sheet2 the active (selected) tab, matching the saved workbook state
$ws2.Activate()
